# UnitTest_CalcTable_Landscape_DataParser.xlsx — "UT - Data Collection Types"
# sheet gets extra data-record rows (different collection sizes per record),
# and the active selection moves from E10 to H9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UT - Data Collection Types")
$ws.Activate()

# --- Row 6: keep only the "list" pair (E6/F6); drop the array/set pair that
#     used to sit alongside it (array now ends at row 9, set now ends at row 9).
$ws.Range("C6:D6").Clear()
$ws.Range("G6:H6").Clear()

# --- Row 7: second data record (#2) — full trio (array/list/set), cloned from
#     row 4's layout, with the sequence number bumped to 2.
$ws.Range("B4:H4").Copy()
$ws.Range("B7:H7").PasteSpecial(-4104)   # xlPasteAll
$ws.Range("B7").Value2 = 2

# --- Row 8: continuation of record #2 (false/BCD), cloned from row 5's layout.
$ws.Range("C5:H5").Copy()
$ws.Range("C8:H8").PasteSpecial(-4104)   # xlPasteAll (gets values/types right)
# Row 5 alternates styles 8/7/8/7/8/7 across C:H; a plain PasteAll into row 8
# collapses them all to style 7, so re-stamp the "false" cells' formats only.
$ws.Range("C5").Copy()
$ws.Range("C8").PasteSpecial(-4122)      # xlPasteFormats
$ws.Range("E5").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("G5").Copy()
$ws.Range("G8").PasteSpecial(-4122)

# --- Row 9: third element of the array (C9/D9) and set (G9/H9) collections —
#     the list collection (E/F) stops at 2 elements, so no E9/F9 here.
$ws.Range("D9").Clear()
$ws.Range("C4:D4").Copy()
$ws.Range("C9:D9").PasteSpecial(-4104)   # xlPasteAll
$ws.Range("G4:H4").Copy()
$ws.Range("G9:H9").PasteSpecial(-4104)   # xlPasteAll

$ws.Range("H9").Select() | Out-Null

$excel.CutCopyMode = 0
